$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is the same value (45171)
# for every data row (2..420). Update it to 45172 for all data rows.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 420 }

$rng = $ws.Range("C2:C$lastRow")
$rng.Value2 = 45172
